{"js": "// Update the date line and the 25 \"NN\u00f7N=\" division prompts in the table,\n// following the same left-to-right, top-to-bottom document order as the\n// non-empty paragraphs in the body (the date paragraph, then each filled\n// table row of five cells).\nconst replacements = [\n  \"2025-02-17 Monday\",\n  \"25\u00f74=\", \"21\u00f77=\", \"98\u00f73=\", \"29\u00f75=\", \"66\u00f77=\",\n  \"97\u00f78=\", \"75\u00f74=\", \"40\u00f77=\", \"84\u00f74=\", \"42\u00f74=\",\n  \"77\u00f74=\", \"97\u00f72=\", \"68\u00f74=\", \"71\u00f72=\", \"28\u00f76=\",\n  \"76\u00f77=\", \"65\u00f79=\", \"22\u00f74=\", \"70\u00f74=\", \"14\u00f77=\",\n  \"29\u00f72=\", \"48\u00f75=\", \"62\u00f79=\", \"79\u00f72=\", \"26\u00f74=\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Keep only paragraphs that actually contain text (skips the blank\n// spacer rows in the table) and apply the replacements in document order.\nconst nonEmpty = paragraphs.items.filter((p) => p.text && p.text.length > 0);\n\nfor (let i = 0; i < replacements.length && i < nonEmpty.length; i++) {\n  nonEmpty[i].insertText(replacements[i], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 \"NN\u00f7N=\" division prompts in the table,\n# following the same left-to-right, top-to-bottom document order as the\n# non-empty paragraphs in the body (the date paragraph, then each filled\n# table row of five cells).\n$replacements = @(\n    \"2025-02-17 Monday\",\n    \"25\u00f74=\", \"21\u00f77=\", \"98\u00f73=\", \"29\u00f75=\", \"66\u00f77=\",\n    \"97\u00f78=\", \"75\u00f74=\", \"40\u00f77=\", \"84\u00f74=\", \"42\u00f74=\",\n    \"77\u00f74=\", \"97\u00f72=\", \"68\u00f74=\", \"71\u00f72=\", \"28\u00f76=\",\n    \"76\u00f77=\", \"65\u00f79=\", \"22\u00f74=\", \"70\u00f74=\", \"14\u00f77=\",\n    \"29\u00f72=\", \"48\u00f75=\", \"62\u00f79=\", \"79\u00f72=\", \"26\u00f74=\"\n)\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n$idx = 0\n\nfor ($i = 1; $i -le $count; $i++) {\n    if ($idx -ge $replacements.Count) { break }\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    $clean = $t -replace \"[\\r\\x07]+$\", \"\"\n    if ($clean.Length -gt 0) {\n        $p.Range.Text = $replacements[$idx]\n        $idx++\n    }\n}\n"}
